$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds the old, not-yet-finished test case
# "FourOfAKind_44442_88885_SecondWin_KeyCard_8" (Joe wins).
# Replace it with the finished test case
# "FourOfAKind_88885_44442_FirstWin_KeyCard_8" (Tom wins) by swapping the
# Tom/Joe hand columns (D/E <-> F/G) and updating the test-case name and
# result text.

# Copy the "black" (non-winning) number format from D2 onto F4, then copy
# the "red" (winning) number format from F5 onto D4, before writing the new
# values - this keeps the existing font table untouched (no new fonts/styles
# get created), matching how Excel itself would recolor the swapped cells.
$ws.Range("D2").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("F5").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("C4").Value = "FourOfAKind_88885_44442_FirstWin_KeyCard_8"
$ws.Range("D4").Value = "8D,8C,8H,8S,5C"
$ws.Range("E4").Value = 88885
$ws.Range("F4").Value = "4D,4S,2D,4H,4C"
$ws.Range("G4").Value = 44442
$ws.Range("H4").Value = "Tom wins. - with four of a kind, key card 8"

$ws.Range("C4").Select()
